# Apply odds updates to Sheet1 per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Cell = "G4"; Value = 1.8 },
    @{ Cell = "P4"; Value = 1.96 },
    @{ Cell = "W4"; Value = 2.24 },
    @{ Cell = "N5"; Value = 3 },
    @{ Cell = "G6"; Value = 1.26 },
    @{ Cell = "I6"; Value = 18 },
    @{ Cell = "L6"; Value = 1.01 },
    @{ Cell = "M6"; Value = 1.02 },
    @{ Cell = "N6"; Value = 3 },
    @{ Cell = "O6"; Value = 1.11 },
    @{ Cell = "R6"; Value = 1.7 },
    @{ Cell = "S6"; Value = 1.84 },
    @{ Cell = "T6"; Value = 1.81 },
    @{ Cell = "U6"; Value = 1.72 },
    @{ Cell = "V6"; Value = 1.05 },
    @{ Cell = "W6"; Value = 4.7 },
    @{ Cell = "X6"; Value = 1000 },
    @{ Cell = "Y6"; Value = 75 },
    @{ Cell = "Z6"; Value = 1000 },
    @{ Cell = "AA6"; Value = 1000 },
    @{ Cell = "AB6"; Value = 18 },
    @{ Cell = "AC6"; Value = 25 },
    @{ Cell = "AD6"; Value = 70 },
    @{ Cell = "AE6"; Value = 1000 },
    @{ Cell = "AF6"; Value = 13.5 },
    @{ Cell = "AG6"; Value = 17.5 },
    @{ Cell = "AH6"; Value = 46 },
    @{ Cell = "AI6"; Value = 1000 },
    @{ Cell = "AJ6"; Value = 13.5 },
    @{ Cell = "AK6"; Value = 19 },
    @{ Cell = "AL6"; Value = 50 },
    @{ Cell = "AM6"; Value = 1000 },
    @{ Cell = "AN6"; Value = 1000 },
    @{ Cell = "AO6"; Value = 1000 },
    @{ Cell = "F7"; Value = 15 },
    @{ Cell = "G7"; Value = 24 },
    @{ Cell = "I7"; Value = 1.23 },
    @{ Cell = "J7"; Value = 7.8 },
    @{ Cell = "L7"; Value = 1.01 },
    @{ Cell = "M7"; Value = 1.02 },
    @{ Cell = "N7"; Value = 7 },
    @{ Cell = "O7"; Value = 1.13 },
    @{ Cell = "R7"; Value = 1.84 },
    @{ Cell = "S7"; Value = 1.9 },
    @{ Cell = "T7"; Value = 2 },
    @{ Cell = "U7"; Value = 1.66 },
    @{ Cell = "V7"; Value = 5.3 },
    @{ Cell = "W7"; Value = 1.04 },
    @{ Cell = "X7"; Value = 48 },
    @{ Cell = "Y7"; Value = 14 },
    @{ Cell = "Z7"; Value = 10.5 },
    @{ Cell = "AA7"; Value = 10.5 },
    @{ Cell = "AB7"; Value = 65 },
    @{ Cell = "AC7"; Value = 25 },
    @{ Cell = "AD7"; Value = 15 },
    @{ Cell = "AE7"; Value = 16.5 },
    @{ Cell = "AF7"; Value = 1000 },
    @{ Cell = "AG7"; Value = 75 },
    @{ Cell = "AH7"; Value = 48 },
    @{ Cell = "AI7"; Value = 48 },
    @{ Cell = "AJ7"; Value = 1000 },
    @{ Cell = "AK7"; Value = 1000 },
    @{ Cell = "AL7"; Value = 1000 },
    @{ Cell = "AM7"; Value = 1000 },
    @{ Cell = "AN7"; Value = 1000 },
    @{ Cell = "AO7"; Value = 3.5 },
    @{ Cell = "F8"; Value = 2.5 },
    @{ Cell = "I8"; Value = 4.1 },
    @{ Cell = "L8"; Value = 1.01 },
    @{ Cell = "M8"; Value = 1.14 },
    @{ Cell = "N8"; Value = 2.16 },
    @{ Cell = "O8"; Value = 1.73 },
    @{ Cell = "Q8"; Value = 3.25 },
    @{ Cell = "R8"; Value = 1.12 },
    @{ Cell = "S8"; Value = 6.4 },
    @{ Cell = "T8"; Value = 2.16 },
    @{ Cell = "U8"; Value = 1.5 },
    @{ Cell = "V8"; Value = 1.32 },
    @{ Cell = "W8"; Value = 1.58 },
    @{ Cell = "X8"; Value = 6.6 },
    @{ Cell = "Y8"; Value = 9 },
    @{ Cell = "Z8"; Value = 25 },
    @{ Cell = "AA8"; Value = 110 },
    @{ Cell = "AB8"; Value = 6.8 },
    @{ Cell = "AC8"; Value = 7.2 },
    @{ Cell = "AD8"; Value = 19 },
    @{ Cell = "AE8"; Value = 1000 },
    @{ Cell = "AF8"; Value = 15 },
    @{ Cell = "AG8"; Value = 14.5 },
    @{ Cell = "AH8"; Value = 30 },
    @{ Cell = "AI8"; Value = 1000 },
    @{ Cell = "AJ8"; Value = 1000 },
    @{ Cell = "AK8"; Value = 1000 },
    @{ Cell = "AL8"; Value = 120 },
    @{ Cell = "AM8"; Value = 1000 },
    @{ Cell = "AN8"; Value = 1000 },
    @{ Cell = "AO8"; Value = 1000 },
    @{ Cell = "H9"; Value = 5.5 },
    @{ Cell = "I9"; Value = 7 },
    @{ Cell = "J9"; Value = 3.2 },
    @{ Cell = "K9"; Value = 3.75 },
    @{ Cell = "L9"; Value = 1.01 },
    @{ Cell = "M9"; Value = 1.09 },
    @{ Cell = "N9"; Value = 2.76 },
    @{ Cell = "P9"; Value = 1.59 },
    @{ Cell = "Q9"; Value = 2.4 },
    @{ Cell = "R9"; Value = 1.22 },
    @{ Cell = "S9"; Value = 4.2 },
    @{ Cell = "T9"; Value = 2.2 },
    @{ Cell = "U9"; Value = 1.71 },
    @{ Cell = "V9"; Value = 1.17 },
    @{ Cell = "W9"; Value = 2.16 },
    @{ Cell = "X9"; Value = 12 },
    @{ Cell = "Y9"; Value = 1000 },
    @{ Cell = "Z9"; Value = 1000 },
    @{ Cell = "AA9"; Value = 1000 },
    @{ Cell = "AB9"; Value = 1000 },
    @{ Cell = "AC9"; Value = 1000 },
    @{ Cell = "AD9"; Value = 1000 },
    @{ Cell = "AE9"; Value = 1000 },
    @{ Cell = "AF9"; Value = 1000 },
    @{ Cell = "AG9"; Value = 1000 },
    @{ Cell = "AH9"; Value = 1000 },
    @{ Cell = "AI9"; Value = 1000 },
    @{ Cell = "AJ9"; Value = 1000 },
    @{ Cell = "AK9"; Value = 1000 },
    @{ Cell = "AL9"; Value = 1000 },
    @{ Cell = "AM9"; Value = 1000 },
    @{ Cell = "AN9"; Value = 1000 },
    @{ Cell = "AO9"; Value = 1000 },
    @{ Cell = "G10"; Value = 1.91 },
    @{ Cell = "J10"; Value = 3.5 },
    @{ Cell = "L10"; Value = 1.01 },
    @{ Cell = "M10"; Value = 1.07 },
    @{ Cell = "N10"; Value = 1.73 },
    @{ Cell = "O10"; Value = 1.38 },
    @{ Cell = "R10"; Value = 1.21 },
    @{ Cell = "S10"; Value = 3.4 },
    @{ Cell = "T10"; Value = 1.01 },
    @{ Cell = "U10"; Value = 1.01 },
    @{ Cell = "V10"; Value = 1.16 },
    @{ Cell = "W10"; Value = 2.12 },
    @{ Cell = "X10"; Value = 1000 },
    @{ Cell = "Y10"; Value = 23 },
    @{ Cell = "Z10"; Value = 65 },
    @{ Cell = "AA10"; Value = 1000 },
    @{ Cell = "AB10"; Value = 10 },
    @{ Cell = "AC10"; Value = 11.5 },
    @{ Cell = "AD10"; Value = 30 },
    @{ Cell = "AE10"; Value = 1000 },
    @{ Cell = "AF10"; Value = 13.5 },
    @{ Cell = "AG10"; Value = 14.5 },
    @{ Cell = "AH10"; Value = 32 },
    @{ Cell = "AI10"; Value = 1000 },
    @{ Cell = "AJ10"; Value = 26 },
    @{ Cell = "AK10"; Value = 29 },
    @{ Cell = "AL10"; Value = 65 },
    @{ Cell = "AM10"; Value = 1000 },
    @{ Cell = "AN10"; Value = 1000 },
    @{ Cell = "AO10"; Value = 1000 },
    @{ Cell = "I11"; Value = 2.82 },
    @{ Cell = "J11"; Value = 2.74 },
    @{ Cell = "F12"; Value = 2.2 },
    @{ Cell = "K12"; Value = 3.6 },
    @{ Cell = "Q12"; Value = 2.08 }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
